$wb = $excel.ActiveWorkbook

# Source sheet that holds the table data to be duplicated
$src = $wb.Worksheets.Item("Tbl_HoldTypesToBeConsideredFor1")

# Add a brand-new worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Tbl_HoldTypesToBeConsideredFor2"

# Copy the table contents (ID / HoldTypesToBeConsideredForCreditCheckReleases)
# from the existing sheet into the newly created sheet
$src.Range("A1:B5").Copy()
$newSheet.Range("A1").PasteSpecial()

# Re-point the defined name used by the table at the new sheet
$definedName = $wb.Names.Item("Tbl_HoldTypesToBeConsideredForCreditCheckReleases")
$definedName.RefersTo = "=Tbl_HoldTypesToBeConsideredFor2!`$A`$1:`$B`$5"

# Make the new sheet the active / selected sheet, matching tabSelected="1"
$newSheet.Select()
